$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (batsman), shifting existing
# D:I columns to F:K.
$ws.Range("D:E").Insert()

# New header cells
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# New data cells
$ws.Range("D2").Value = "Sunrisers Hyderabad"
$ws.Range("E2").Value = "Kolkata Knight Riders"
